$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data scraped on Tue Dec  5 05:28:07 UTC 2023
$ws.Range("D2").Value = "41.876.29"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "2.230.49"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'232.93"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").Value = "'0.620"
$ws.Range("E6").Value = "  -2.72%  "
$ws.Range("D7").Value = "'60.71"
$ws.Range("E7").Value = "  -6.98%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.406"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").Value = "'58.20"
$ws.Range("E10").Value = "  -4.23%  "
$ws.Range("D11").Value = "'0.0903"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "2.560.65"
$ws.Range("E13").Value = "  -0.76%  "
$ws.Range("D14").Value = "'15.55"
$ws.Range("E14").Value = "  -4.03%  "
$ws.Range("D15").Value = "'22.77"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "'5.65"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("E17").Value = "  -3.65%  "
$ws.Range("D18").Value = "2.242.93"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "41.807.18"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").Value = "0.0₃0912"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").Value = "'72.66"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("D22").Value = "'6.19"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "'248.32"
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "'2.39"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").Value = "'2.32"
$ws.Range("E26").Value = "  -3.55%  "
$ws.Range("D27").Value = "'9.70"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").Value = "'169.72"
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("E29").Value = "  -3.04%  "
$ws.Range("D30").Value = "'19.94"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D32").Value = "'2.58"
$ws.Range("E32").Value = "  -9.65%  "
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("D34").Value = "'5.03"
$ws.Range("E34").Value = "  +2.77%  "
$ws.Range("D35").Value = "'4.72"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("D36").Value = "'0.0660"
$ws.Range("E36").Value = "  +4.04%  "
$ws.Range("D37").Value = "'6.59"
$ws.Range("E37").Value = "  -9.42%  "
$ws.Range("D38").Value = "'2.41"
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("E39").Value = "  -5.80%  "
$ws.Range("B40").Value = "TerraClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D40").Value = "'0.000248"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("B41").Value = "BinanceUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D41").Value = "'0.997"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "'0.0241"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").Value = "'8.66"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'1.23"
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.50"
$ws.Range("E45").Value = "  -8.21%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'99.16"
$ws.Range("E46").Value = "  -3.25%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0969"
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("D48").Value = "1.471.88"
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("D49").Value = "'16.68"
$ws.Range("E49").Value = "  -5.67%  "
$ws.Range("D50").Value = "'2.79"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").Value = "'2.28"
$ws.Range("E51").Value = "  +7.13%  "
